$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing rows down
$ws.Rows.Item(2).Insert()

# Set the new cell's content
$ws.Cells.Item(2, 1).Value = "objetivos, problema, justificativa e metodologia são subtopicos de introdução"
